# Update the Price (D) and Volume(1h) (E) columns of the cryptos list
# with the latest scraped values. Numeric-looking Price values are
# prefixed with a leading apostrophe so Excel keeps them stored as text
# (matching the original sheet layout) instead of auto-converting them
# to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.279.08'
$ws.Range("E2").Value = '  +2.41%  '
$ws.Range("D3").Value = '2.423.11'
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'555.52"
$ws.Range("E5").Value = '  +2.09%  '
$ws.Range("D6").Value = "'143.16"
$ws.Range("E6").Value = '  +4.64%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'0.533"
$ws.Range("E8").Value = '  +1.55%  '
$ws.Range("D9").Value = '2.422.06'
$ws.Range("E9").Value = '  +3.21%  '
$ws.Range("D10").Value = "'0.109"
$ws.Range("E10").Value = '  +4.52%  '
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("E12").Value = '  +1.45%  '
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("D14").Value = "'26.23"
$ws.Range("E14").Value = '  +6.39%  '
$ws.Range("D15").Value = "'0.0000175"
$ws.Range("E15").Value = '  +9.11%  '
$ws.Range("D16").Value = '2.862.03'
$ws.Range("E16").Value = '  +3.16%  '
$ws.Range("D17").Value = '62.209.08'
$ws.Range("E17").Value = '  +2.44%  '
$ws.Range("D18").Value = '2.422.58'
$ws.Range("E18").Value = '  +3.13%  '
$ws.Range("E19").Value = '  +4.29%  '
$ws.Range("E20").Value = '  +1.94%  '
$ws.Range("D21").Value = "'324.39"
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("E22").Value = '  +2.54%  '
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("E24").Value = '  +5.32%  '
$ws.Range("D25").Value = "'64.87"
$ws.Range("E25").Value = '  +2.37%  '
$ws.Range("E26").Value = '  +9.05%  '
$ws.Range("D27").Value = "'573.05"
$ws.Range("E27").Value = '  +16.00%  '
$ws.Range("D28").Value = '2.541.04'
$ws.Range("E28").Value = '  +3.05%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  +5.44%  '
$ws.Range("D31").Value = '0.0₃0934'
$ws.Range("E31").Value = '  +8.97%  '
$ws.Range("E32").Value = '  +6.15%  '
$ws.Range("E33").Value = '  +2.13%  '
$ws.Range("E34").Value = '  +4.13%  '
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = '  +4.04%  '
$ws.Range("D36").Value = "'5.69"
$ws.Range("E36").Value = '  +8.95%  '
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = "'4.82"
$ws.Range("E38").Value = '  +5.01%  '
$ws.Range("D39").Value = "'0.384"
$ws.Range("E39").Value = '  +2.29%  '
$ws.Range("E40").Value = '  +4.25%  '
$ws.Range("D41").Value = "'18.74"
$ws.Range("E41").Value = '  +1.48%  '
$ws.Range("D42").Value = "'149.56"
$ws.Range("E42").Value = '  +3.17%  '
$ws.Range("D44").Value = "'41.69"
$ws.Range("E44").Value = '  +2.64%  '
$ws.Range("D45").Value = "'2.31"
$ws.Range("E45").Value = '  +14.43%  '
$ws.Range("D46").Value = "'150.96"
$ws.Range("E46").Value = '  +5.54%  '
$ws.Range("D47").Value = "'3.63"
$ws.Range("E47").Value = '  +2.09%  '
$ws.Range("E48").Value = '  +4.75%  '
$ws.Range("D49").Value = "'20.37"
$ws.Range("E49").Value = '  +7.01%  '
$ws.Range("E50").Value = '  +3.70%  '
$ws.Range("E51").Value = '  +1.60%  '
